$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.965
$ws.Range("B9").Value = 5.317
$ws.Range("C12").Value = -11.171
$ws.Range("E15").Value = 16.346
$ws.Range("B18").Value = 5.19
$ws.Range("B20").Value = 6.952000000000001
$ws.Range("C26").Value = -13.563
$ws.Range("B27").Value = 6.237
$ws.Range("C27").Value = -13.27
$ws.Range("C29").Value = -12.735
$ws.Range("C37").Value = -13.351
$ws.Range("C38").Value = -12.986
$ws.Range("E38").Value = 16.556
$ws.Range("E44").Value = 16.814
$ws.Range("C51").Value = -11.471
$ws.Range("E51").Value = 17.207
$ws.Range("C55").Value = -13.752
$ws.Range("E57").Value = 16.366
$ws.Range("E63").Value = 17.601
$ws.Range("B69").Value = 5.425
$ws.Range("C69").Value = -11.174
$ws.Range("C70").Value = -12.091
$ws.Range("E70").Value = 17.447
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B82").Value = 5.241
$ws.Range("C83").Value = -13.169
$ws.Range("E99").Value = 16.601
$ws.Range("C102").Value = -13.419
